# Update gh-pages output - apply scraped bilibili concert/exhibition data refresh
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (Exhibitions) - refresh "想去人数" (interested count) column F
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(4,6).Value2 = 582
$ws1.Cells.Item(5,6).Value2 = 9237
$ws1.Cells.Item(7,6).Value2 = 11971
$ws1.Cells.Item(8,6).Value2 = 11971
$ws1.Cells.Item(11,6).Value2 = 39
$ws1.Cells.Item(12,6).Value2 = 125
$ws1.Cells.Item(14,6).Value2 = 445
$ws1.Cells.Item(16,6).Value2 = 2051
$ws1.Cells.Item(17,6).Value2 = 816
$ws1.Cells.Item(18,6).Value2 = 776
$ws1.Cells.Item(19,6).Value2 = 394
$ws1.Cells.Item(20,6).Value2 = 38
$ws1.Cells.Item(21,6).Value2 = 409
$ws1.Cells.Item(24,6).Value2 = 663
$ws1.Cells.Item(26,6).Value2 = 16
$ws1.Cells.Item(27,6).Value2 = 1535
$ws1.Cells.Item(29,6).Value2 = 22
$ws1.Cells.Item(32,6).Value2 = 504
$ws1.Cells.Item(33,6).Value2 = 1443
$ws1.Cells.Item(37,6).Value2 = 526
$ws1.Cells.Item(38,6).Value2 = 395
$ws1.Cells.Item(39,6).Value2 = 2172
$ws1.Cells.Item(40,6).Value2 = 375
$ws1.Cells.Item(41,6).Value2 = 74
$ws1.Cells.Item(43,6).Value2 = 575
$ws1.Cells.Item(44,6).Value2 = 446
$ws1.Cells.Item(45,6).Value2 = 158
$ws1.Cells.Item(46,6).Value2 = 882
$ws1.Cells.Item(49,6).Value2 = 294
$ws1.Cells.Item(50,6).Value2 = 268

# ---------------------------------------------------------------------------
# Sheet "演出" (Performances) - the VGL show (row 3, id=89083) was removed
# from the upstream feed. Shift every later row's content (columns B-I) up
# by one row (column A keeps its original sequential index numbers), then
# delete the now-duplicate last row (26). Afterwards re-apply the "想去人数"
# increments that also occurred between scrapes.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")

# Column B holds dates formatted as plain text (e.g. "2024-10-02"). Force
# text format before the bulk copy so Excel does not reinterpret them as
# real date serials, then drop back to the default "Normal" style so no
# stray formatting is left behind.
$bDst = $ws2.Range("B3:B25")
$bDst.NumberFormat = "@"
$bDst.Value2 = $ws2.Range("B4:B26").Value2
$bDst.Style = "Normal"

# Columns C:I can be copied directly as-is.
$ws2.Range("C3:I25").Value2 = $ws2.Range("C4:I26").Value2

$ws2.Rows.Item(26).Delete()

# Additional "想去人数" increments on top of the shifted data
$ws2.Cells.Item(17,6).Value2 = 85    # majiko中国巡演-2024 (id=92300)
$ws2.Cells.Item(23,6).Value2 = 108   # Ayasa LIVE TOUR 2024 (id=92778)
$ws2.Cells.Item(24,6).Value2 = 66    # 吹响吧！ACG！ (id=93059)
$ws2.Cells.Item(25,6).Value2 = 425   # 花たん 2024 LIVE in Beijing (id=90341)

# ---------------------------------------------------------------------------
# Sheet "本地生活" (Local life) - refresh "想去人数" column F
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(6,6).Value2 = 258

# ---------------------------------------------------------------------------
# Sheet "全部类型" (All types, merged view) - refresh "想去人数" column F
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(7,6).Value2 = 582
$ws4.Cells.Item(8,6).Value2 = 9237
$ws4.Cells.Item(10,6).Value2 = 11971
$ws4.Cells.Item(11,6).Value2 = 39
$ws4.Cells.Item(12,6).Value2 = 125
$ws4.Cells.Item(13,6).Value2 = 445
$ws4.Cells.Item(14,6).Value2 = 2051
$ws4.Cells.Item(15,6).Value2 = 816
$ws4.Cells.Item(16,6).Value2 = 776
$ws4.Cells.Item(17,6).Value2 = 394
$ws4.Cells.Item(18,6).Value2 = 38
$ws4.Cells.Item(19,6).Value2 = 409
$ws4.Cells.Item(22,6).Value2 = 663
$ws4.Cells.Item(25,6).Value2 = 16
$ws4.Cells.Item(26,6).Value2 = 258
$ws4.Cells.Item(27,6).Value2 = 1535
$ws4.Cells.Item(30,6).Value2 = 504
$ws4.Cells.Item(32,6).Value2 = 1443
$ws4.Cells.Item(36,6).Value2 = 85
$ws4.Cells.Item(37,6).Value2 = 526
$ws4.Cells.Item(38,6).Value2 = 395
$ws4.Cells.Item(39,6).Value2 = 2174
$ws4.Cells.Item(40,6).Value2 = 74
$ws4.Cells.Item(42,6).Value2 = 575
$ws4.Cells.Item(43,6).Value2 = 446
$ws4.Cells.Item(44,6).Value2 = 158
$ws4.Cells.Item(46,6).Value2 = 66
$ws4.Cells.Item(47,6).Value2 = 425
$ws4.Cells.Item(49,6).Value2 = 294
$ws4.Cells.Item(50,6).Value2 = 268
